$wb = $excel.ActiveWorkbook

# Rename sheet "getDataGraphQL" to "getDataGraphQL-bc"
$wsGetData = $wb.Worksheets.Item("getDataGraphQL")
$wsGetData.Name = "getDataGraphQL-bc"

# Activate it and select B29 (it becomes the active/tabSelected sheet)
$wsGetData.Activate()
$wsGetData.Range("B29").Select()
